$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows for line7 / line8, right after the existing line6 row (row 7),
# pushing the extr1..extr8 rows down from 8-15 to 10-17.
$ws.Rows.Item(8).Resize(2).Insert()

# New row 8: line7
$ws.Cells.Item(8, 1).Value = 6
$ws.Cells.Item(8, 2).Value = "line7"
$ws.Cells.Item(8, 3).Value = 14
$ws.Cells.Item(8, 4).Value = 11
$ws.Cells.Item(8, 5).Value = $true

# New row 9: line8
$ws.Cells.Item(9, 1).Value = 7
$ws.Cells.Item(9, 2).Value = "line8"
$ws.Cells.Item(9, 3).Value = 16
$ws.Cells.Item(9, 4).Value = 9
$ws.Cells.Item(9, 5).Value = $true

# Copy style (bold/border/center) used for column-A id cells onto the new rows' A cells.
$ws.Cells.Item(7, 1).Copy()
$ws.Cells.Item(8, 1).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(9, 1).PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Refresh the id (A) / from_bus (C) / to_bus (D) / in_service (E) values for the
# extr1..extr8 rows, which shifted down from rows 8-15 to 10-17.
$ws.Cells.Item(10, 1).Value = 8
$ws.Cells.Item(10, 3).Value = 5
$ws.Cells.Item(10, 4).Value = 12
$ws.Cells.Item(10, 5).Value = $true   # extr1

$ws.Cells.Item(11, 1).Value = 9
$ws.Cells.Item(11, 3).Value = 5
$ws.Cells.Item(11, 4).Value = 9
$ws.Cells.Item(11, 5).Value = $true   # extr2

$ws.Cells.Item(12, 1).Value = 10
$ws.Cells.Item(12, 3).Value = 10
$ws.Cells.Item(12, 4).Value = 11
$ws.Cells.Item(12, 5).Value = $false  # extr3

$ws.Cells.Item(13, 1).Value = 11
$ws.Cells.Item(13, 3).Value = 7
$ws.Cells.Item(13, 4).Value = 8
$ws.Cells.Item(13, 5).Value = $false  # extr4

$ws.Cells.Item(14, 1).Value = 12
$ws.Cells.Item(14, 3).Value = 9
$ws.Cells.Item(14, 4).Value = 11
$ws.Cells.Item(14, 5).Value = $false  # extr5

$ws.Cells.Item(15, 1).Value = 13
$ws.Cells.Item(15, 3).Value = 7
$ws.Cells.Item(15, 4).Value = 11
$ws.Cells.Item(15, 5).Value = $true   # extr6

$ws.Cells.Item(16, 1).Value = 14
$ws.Cells.Item(16, 3).Value = 5
$ws.Cells.Item(16, 4).Value = 7
$ws.Cells.Item(16, 5).Value = $true   # extr7

$ws.Cells.Item(17, 1).Value = 15
$ws.Cells.Item(17, 3).Value = 8
$ws.Cells.Item(17, 4).Value = 5
$ws.Cells.Item(17, 5).Value = $false  # extr8
